$d = $word.ActiveDocument

# --- Paragraph 2 rewrite -------------------------------------------------

$d.Content.Find.Execute("Over the last 2 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "I wish to draw your attention to a ", 2) | Out-Null

$d.Content.Find.Execute("years", $true, $false, $false, $false, $false,
                         $true, 1, $false, "problem which", 2) | Out-Null

$d.Content.Find.Execute(" I have purchased several quality books through your website, and have been very pleased with the quality and service I have received. On September 14, 2018, I used your website to order a first edition copy of Service for Two by Madeleine ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " arose due to book ‘Service for Two’ by Madeleine ", 2) | Out-Null

$d.Content.Find.Execute(", published in 1964 by Carlton Publishers. My online confirmation number for the order is F123U456789. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " I ordered on your online shop on 28th September. ", 2) | Out-Null

# --- Insert the _GoBack bookmark where the edit actually happened -------
# (between "...problem which a" and "rose due to book...")

$bmRange = $d.Content
$bmRange.Find.Execute("problem which a", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0) | Out-Null
$bmRange.Start = $bmRange.End
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- Paragraph 3: merge "On " / "September " into a single run ----------

$d.Content.Find.Execute("On September ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "On September ", 2) | Out-Null

# --- Remove the old trailing _GoBack bookmark next to the signature -----
# (Bookmarks.Add above already re-defines "_GoBack" at the new location,
#  which automatically removes the previous one sharing that name.)

Write-Output $d.Content.Text
